# Update crypto price/volume figures per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores figures as text (e.g. "59.610.18", "0.339"),
# not numbers. Pre-format the numeric-looking ones as Text so Excel
# does not auto-convert them (and silently drop trailing zeros, e.g.
# "0.340" -> 0.34).
$textCells = @("D5", "D6", "D12", "D13", "D18", "D20", "D23", "D26", "D27", "D29", "D31", "D34", "D40", "D42", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '59.539.19'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '2.340.00'
$ws.Range("E3").Value = '  -3.29%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '558.07'
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("D6").Value = '132.23'
$ws.Range("E6").Value = '  -3.69%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -3.13%  '
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("D12").Value = '0.340'
$ws.Range("E12").Value = '  -4.03%  '
$ws.Range("D13").Value = '23.88'
$ws.Range("E13").Value = '  -6.15%  '
$ws.Range("D14").Value = '2.763.27'
$ws.Range("E14").Value = '  -2.92%  '
$ws.Range("D15").Value = '59.517.32'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").Value = '2.343.46'
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").Value = '10.95'
$ws.Range("E18").Value = '  -3.25%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '318.59'
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '63.82'
$ws.Range("E23").Value = '  -3.01%  '
$ws.Range("E24").Value = '  -2.99%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '8.36'
$ws.Range("D27").Value = '1.35'
$ws.Range("E27").Value = '  -2.04%  '
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").Value = '170.99'
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("E30").Value = '  -4.06%  '
$ws.Range("D31").Value = '5.95'
$ws.Range("E31").Value = '  -1.84%  '
$ws.Range("E32").Value = '  +6.49%  '
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("D34").Value = '17.97'
$ws.Range("E34").Value = '  -3.39%  '
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -3.26%  '
$ws.Range("E39").Value = '  -2.17%  '
$ws.Range("D40").Value = '313.73'
$ws.Range("E40").Value = '  -2.51%  '
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").Value = '144.37'
$ws.Range("E42").Value = '  +2.98%  '
$ws.Range("E43").Value = '  -5.31%  '
$ws.Range("D44").Value = '0.0957'
$ws.Range("E44").Value = '  -1.53%  '
$ws.Range("D45").Value = '0.0505'
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("D46").Value = '18.95'
$ws.Range("E46").Value = '  -3.42%  '
$ws.Range("D47").Value = '0.561'
$ws.Range("E47").Value = '  -3.20%  '
$ws.Range("E48").Value = '  -3.25%  '
$ws.Range("D49").Value = '11.07'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '0.940'
$ws.Range("E51").Value = '  -0.28%  '
